# Add validation check before vlan creation - cisco & junos
# Populate the "group" column for the Cisco/Junos devices that were
# previously left as "none", and correct the "type" column for two
# rows that were mis-tagged as "backbone" instead of "access".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "access"

$ws.Range("B7").Value = "junos"
$ws.Range("E7").Value = "access"

$ws.Range("B8").Value = "cisco"

# Restore the active selection recorded in the sheet view.
$ws.Range("K9").Select()
